$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 4542.2856
$ws.Range("I8").Value = 4542.2856
$ws.Range("K8").Value = 13626.8568
$ws.Range("M8").Value = -13487.8568
$ws.Range("H17").Value = 1727.4546
$ws.Range("J17").Value = 1727.4546
$ws.Range("L17").Value = 5182.3638
$ws.Range("N17").Value = -5518.3638
$ws.Range("H46").Value = 1730.1
$ws.Range("I46").Value = 1200.25
$ws.Range("J46").Value = 2083.3333
$ws.Range("K46").Value = 3600.75
$ws.Range("L46").Value = 6249.999899999999
$ws.Range("M46").Value = -3481.75
$ws.Range("N46").Value = -6487.999899999999
$ws.Range("H60").Value = 1730.1
$ws.Range("I60").Value = 1200.25
$ws.Range("J60").Value = 2083.3333
$ws.Range("K60").Value = 3600.75
$ws.Range("L60").Value = 6249.999899999999
$ws.Range("M60").Value = -3116.75
$ws.Range("N60").Value = -7217.999899999999
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H135").Value = 1909.1
$ws.Range("I135").Value = 2248.6667
$ws.Range("J135").Value = 1399.75
$ws.Range("K135").Value = 20238.0003
$ws.Range("L135").Value = 12597.75
$ws.Range("M135").Value = -17703.0003
$ws.Range("N135").Value = -17667.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 715.1429000000001
$ws.Range("I2").Value = 875
$ws.Range("J2").Value = 502
$ws.Range("K2").Value = 875
$ws.Range("L2").Value = 502
$ws.Range("M2").Value = -762
$ws.Range("N2").Value = -728
$ws.Range("H32").Value = 10879854
$ws.Range("I32").Value = 17861534
$ws.Range("J32").Value = 19464.334
$ws.Range("K32").Value = 17861534
$ws.Range("L32").Value = 19464.334
$ws.Range("M32").Value = -17861247
$ws.Range("N32").Value = -20038.334
$ws.Range("H74").Value = 8937284
$ws.Range("I74").Value = 13891419
$ws.Range("K74").Value = 13891419
$ws.Range("M74").Value = -13890545
$ws.Range("H77").Value = 8937284
$ws.Range("I77").Value = 13891419
$ws.Range("K77").Value = 69457095
$ws.Range("M77").Value = -69452727
$ws.Range("H116").Value = 715.1429000000001
$ws.Range("I116").Value = 875
$ws.Range("J116").Value = 502
$ws.Range("K116").Value = 875
$ws.Range("L116").Value = 502
$ws.Range("M116").Value = 1419
$ws.Range("N116").Value = -5090
$ws.Range("H132").Value = 4671.316
$ws.Range("I132").Value = 2213.4688
$ws.Range("K132").Value = 6640.4064
$ws.Range("M132").Value = -4110.4064

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 715.1429000000001
$ws.Range("I3").Value = 875
$ws.Range("J3").Value = 502
$ws.Range("K3").Value = 875
$ws.Range("L3").Value = 502
$ws.Range("M3").Value = -761
$ws.Range("N3").Value = -730
$ws.Range("H36").Value = 3014.1667
$ws.Range("I36").Value = 2376.6667
$ws.Range("J36").Value = 3651.6667
$ws.Range("K36").Value = 2376.6667
$ws.Range("L36").Value = 3651.6667
$ws.Range("M36").Value = -1842.6667
$ws.Range("N36").Value = -4719.6667
$ws.Range("H117").Value = 115000
$ws.Range("J117").Value = 115000
$ws.Range("L117").Value = 115000
$ws.Range("N117").Value = -124178
$ws.Range("H134").Value = 206038.48
$ws.Range("J134").Value = 851887
$ws.Range("L134").Value = 2555661
$ws.Range("N134").Value = -2560731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 593844.5600000001
$ws.Range("J31").Value = 1178653.2
$ws.Range("L31").Value = 1178653.2
$ws.Range("N31").Value = -1179243.2
$ws.Range("H34").Value = 593844.5600000001
$ws.Range("J34").Value = 1178653.2
$ws.Range("L34").Value = 1178653.2
$ws.Range("N34").Value = -1179057.2
$ws.Range("H105").Value = 2338.889
$ws.Range("I105").Value = 2405.2856
$ws.Range("J105").Value = 2106.5
$ws.Range("K105").Value = 2405.2856
$ws.Range("L105").Value = 2106.5
$ws.Range("M105").Value = -658.2856000000002
$ws.Range("N105").Value = -5600.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 8951.556
$ws.Range("I45").Value = 26000.5
$ws.Range("J45").Value = 4080.4285
$ws.Range("K45").Value = 78001.5
$ws.Range("L45").Value = 12241.2855
$ws.Range("M45").Value = -77469.5
$ws.Range("N45").Value = -13305.2855
$ws.Range("H113").Value = 1323.6471
$ws.Range("I113").Value = 759.4
$ws.Range("J113").Value = 1558.75
$ws.Range("K113").Value = 2278.2
$ws.Range("L113").Value = 4676.25
$ws.Range("M113").Value = -108.1999999999998
$ws.Range("N113").Value = -9016.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 101666.664
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224
$ws.Range("H119").Value = 111000
$ws.Range("J119").Value = 111000
$ws.Range("L119").Value = 111000
$ws.Range("N119").Value = -120676
$ws.Range("H130").Value = 548500
$ws.Range("J130").Value = 1012000
$ws.Range("L130").Value = 1012000
$ws.Range("N130").Value = -1022040
$ws.Range("H132").Value = 27780116
$ws.Range("I132").Value = 29414062
$ws.Range("K132").Value = 88242186
$ws.Range("M132").Value = -88239656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15637.389
$ws.Range("I7").Value = 12447.5
$ws.Range("J7").Value = 18189.3
$ws.Range("K7").Value = 12447.5
$ws.Range("L7").Value = 18189.3
$ws.Range("M7").Value = -12335.5
$ws.Range("N7").Value = -18413.3
$ws.Range("H10").Value = 18333
$ws.Range("J10").Value = 18333
$ws.Range("L10").Value = 18333
$ws.Range("N10").Value = -18613
$ws.Range("H46").Value = 2819.3333
$ws.Range("I46").Value = 2175.25
$ws.Range("K46").Value = 2175.25
$ws.Range("M46").Value = -1987.25
$ws.Range("H55").Value = 111112056
$ws.Range("I55").Value = 142858080
$ws.Range("J55").Value = 1001
$ws.Range("K55").Value = 142858080
$ws.Range("L55").Value = 1001
$ws.Range("M55").Value = -142857907
$ws.Range("N55").Value = -1347
$ws.Range("H61").Value = 1080.3448
$ws.Range("I61").Value = 825.1053000000001
$ws.Range("J61").Value = 1565.3
$ws.Range("K61").Value = 825.1053000000001
$ws.Range("L61").Value = 1565.3
$ws.Range("M61").Value = -623.1053000000001
$ws.Range("N61").Value = -1969.3
$ws.Range("H68").Value = 3689
$ws.Range("I68").Value = 3748.6
$ws.Range("J68").Value = 3518.7144
$ws.Range("K68").Value = 3748.6
$ws.Range("L68").Value = 3518.7144
$ws.Range("M68").Value = -2999.6
$ws.Range("N68").Value = -5016.7144
$ws.Range("H71").Value = 3689
$ws.Range("I71").Value = 3748.6
$ws.Range("J71").Value = 3518.7144
$ws.Range("K71").Value = 18743
$ws.Range("L71").Value = 17593.572
$ws.Range("M71").Value = -14999
$ws.Range("N71").Value = -25081.572
$ws.Range("H87").Value = 441272.66
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 441272.66
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 441272.66
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -443518.66
$ws.Range("H90").Value = 441272.66
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 441272.66
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 1323817.98
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -1335049.98
$ws.Range("H104").Value = 36934.75
$ws.Range("J104").Value = 36934.75
$ws.Range("L104").Value = 36934.75
$ws.Range("N104").Value = -43922.75
$ws.Range("H113").Value = 1080.3448
$ws.Range("I113").Value = 825.1053000000001
$ws.Range("J113").Value = 1565.3
$ws.Range("K113").Value = 825.1053000000001
$ws.Range("L113").Value = 1565.3
$ws.Range("M113").Value = 1344.8947
$ws.Range("N113").Value = -5905.3
$ws.Range("H126").Value = 15637.389
$ws.Range("I126").Value = 12447.5
$ws.Range("J126").Value = 18189.3
$ws.Range("K126").Value = 37342.5
$ws.Range("L126").Value = 54567.89999999999
$ws.Range("M126").Value = -34872.5
$ws.Range("N126").Value = -59507.89999999999
$ws.Range("H132").Value = 407638.2
$ws.Range("I132").Value = 8712
$ws.Range("K132").Value = 26136
$ws.Range("M132").Value = -23606
$ws.Range("H136").Value = 124316.7
$ws.Range("I136").Value = 2698.8333
$ws.Range("K136").Value = 8096.499899999999
$ws.Range("M136").Value = -5546.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 23572.25
$ws.Range("I58").Value = 16598
$ws.Range("J58").Value = 44495
$ws.Range("K58").Value = 16598
$ws.Range("L58").Value = 44495
$ws.Range("M58").Value = -16290
$ws.Range("N58").Value = -45111
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H132").Value = 307073.1
$ws.Range("I132").Value = 3597
$ws.Range("K132").Value = 10791
$ws.Range("M132").Value = -8261
